$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A25").Value = "UBX-40"
$ws.Range("B25").Value = "UBX-40 USRP DAUGHTERBOARD (10 MHZ - 6 GHZ, 40 MHZ BW) - ETTUS RESEARCH"
$ws.Range("C25").Value = "PA1288699"
$ws.Range("E25").Value = "San Diego - Palomar"

$ws.Range("A26").Value = "UBX-40"
$ws.Range("B26").Value = "UBX-40 USRP DAUGHTERBOARD (10 MHZ - 6 GHZ, 40 MHZ BW) - ETTUS RESEARCH"
$ws.Range("C26").Value = "PA1288706"
$ws.Range("E26").Value = "San Diego - Palomar"

$ws.Range("A25").Font.Color = 0
$ws.Range("E25").Font.Color = 0
$ws.Range("A26").Font.Color = 0
$ws.Range("E26").Font.Color = 0

$ws.Range("C29").Select() | Out-Null
